# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" between "总计" and "2021-Q3" and fill it
#    with the fund-holding breakdown for the new quarter.
# 2. Update the "总计" (totals) sheet: push the existing 2021-Q3 summary row
#    down one row and insert a new summary row for 2022-Q4 above it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert the new 2022-Q4 summary row, shift 2021-Q3 down.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Copy the existing 2021-Q3 row down into row 3 first (style + values), then
# overwrite row 2 with the new 2022-Q4 figures. A column keeps incrementing
# as a 0-based running index, matching the rest of the sheet.
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q3"
$totals.Range("C3").Value = 3
$totals.Range("D3").Value = 0.66

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 7
$totals.Range("D2").Value = 1.27

# ---------------------------------------------------------------------------
# 2. New "2022-Q4" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totals)
$newSheet.Name = "2022-Q4"

# Pull header formatting (bold/centered/bordered) from the totals sheet's
# header cell so the new sheet reuses the existing style instead of creating
# a new one. The "A" index column shares the same style as the header row
# throughout this workbook, so seed A2:A8 with it too.
$totals.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("A2:A8").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

function Set-TextCell($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$rows = @(
    @{A=0; B="513500"; C="博时标普500ETF（QDII）";              D="93.59"; E="95.36"; F="1.17"; G="1.0950"; H=10},
    @{A=1; B="012860"; C="易方达标普500指数（QDII-LOF）人民币 C"; D="4.75";  E="91.65"; F="1.12"; G="0.0532"; H=10},
    @{A=2; B="161125"; C="易方达标普500指数（QDII-LOF）人民币";    D="4.75";  E="91.65"; F="1.12"; G="0.0532"; H=10},
    @{A=3; B="003718"; C="易方达标普500指数（QDII-LOF）美元A";    D="4.65";  E="91.65"; F="1.12"; G="0.0521"; H=10},
    @{A=4; B="159612"; C="国泰标普500ETF（QDII）";              D="0.86";  E="94.21"; F="1.14"; G="0.0098"; H=10},
    @{A=5; B="159655"; C="华夏标普500ETF（QDII）";              D="0.21";  E="93.70"; F="1.15"; G="0.0024"; H=9},
    @{A=6; B="012861"; C="易方达标普500指数（QDII-LOF）美元 C";  D="0.10";  E="91.65"; F="1.12"; G="0.0011"; H=10}
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row.A

    Set-TextCell $newSheet.Range("B$r") $row.B
    Set-TextCell $newSheet.Range("C$r") $row.C
    Set-TextCell $newSheet.Range("D$r") $row.D
    Set-TextCell $newSheet.Range("E$r") $row.E
    Set-TextCell $newSheet.Range("F$r") $row.F
    Set-TextCell $newSheet.Range("G$r") $row.G
    $newSheet.Range("H$r").Value = $row.H

    $r = $r + 1
}

# Restore "2021-Q3" as the selected/active sheet (unchanged from before the
# edit; adding a new sheet shifts the active tab by default). Sheet
# references are positional, so re-fetch by name now that all inserts are
# done rather than reusing a handle captured earlier.
$wb.Worksheets.Item("2021-Q3").Activate()
